$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Header row: add new "Browser"/"OS" labels over existing I1/J1 cells, and
#     create a brand-new K1 "Comments" header (copy formatting from J1 first
#     so the new cell keeps the same style as the rest of the header row). ---
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial($xlPasteFormats)

$ws.Range("I1").Value = "Browser"
$ws.Range("J1").Value = "OS"
$ws.Range("K1").Value = "Comments"

# --- Move the per-row Browser values ("chrome") from column J to column I.
#     Column I has no existing cells on rows 2-8, so copy formatting from the
#     neighboring H cell (style "s=1") before writing the new values. ---
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial($xlPasteFormats)
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial($xlPasteFormats)
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial($xlPasteFormats)
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial($xlPasteFormats)
$ws.Range("H6").Copy()
$ws.Range("I6").PasteSpecial($xlPasteFormats)
$ws.Range("H7").Copy()
$ws.Range("I7").PasteSpecial($xlPasteFormats)
$ws.Range("H8").Copy()
$ws.Range("I8").PasteSpecial($xlPasteFormats)

$ws.Range("I2").Value = "chrome"
$ws.Range("I3").Value = "chrome"
$ws.Range("I4").Value = "chrome"
$ws.Range("I5").Value = "chrome"
$ws.Range("I6").Value = "chrome"
$ws.Range("I7").Value = "chrome"
$ws.Range("I8").Value = "chrome"

# Remove the old column J cells entirely (now the empty "OS" column, no data)
$ws.Range("J2").Clear()
$ws.Range("J3").Clear()
$ws.Range("J4").Clear()
$ws.Range("J5").Clear()
$ws.Range("J6").Clear()
$ws.Range("J7").Clear()
$ws.Range("J8").Clear()

# --- Identifier column (E): the "openurl" step (row 2) no longer needs an identifier ---
$ws.Range("E2").Clear()

# Fill in the missing "xpath" identifiers for rows that previously had none
$ws.Range("E5").Value = "xpath"
$ws.Range("E6").Value = "xpath"
$ws.Range("E7").Value = "xpath"
